$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.321.48'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.41%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.552.25'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.39%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.13%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '209.66'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.63%  '

# Row 6
$ws.Range('E6').Value = '  -1.91%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.09%  '

# Row 8
$ws.Range('E8').Value = '  -2.26%  '

# Row 9
$ws.Range('E9').Value = '  -2.16%  '

# Row 10
$ws.Range('E10').Value = '  -1.57%  '

# Row 11
$ws.Range('E11').Value = '  -0.18%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.773.33'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.39%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.562.35'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.74%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '28.304.41'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.44%  '

# Row 15
$ws.Range('E15').Value = '  -1.69%  '

# Row 16
$ws.Range('E16').Value = '  -2.60%  '

# Row 17
$ws.Range('E17').Value = '  -3.08%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '228.09'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.15%  '

# Row 19
$ws.Range('E19').Value = '  -0.97%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0674'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.83%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.999'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.08%  '

# Row 22
$ws.Range('E22').Value = '  +0.48%  '

# Row 23
$ws.Range('E23').Value = '  -3.27%  '

# Row 24
$ws.Range('E24').Value = '  -4.43%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '150.93'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.68%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '14.76'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.94%  '

# Row 27
$ws.Range('E27').Value = '  -1.19%  '

# Row 28
$ws.Range('E28').Value = '  -0.13%  '

# Row 29
$ws.Range('E29').Value = '  -3.18%  '

# Row 30
$ws.Range('E30').Value = '  -4.04%  '

# Row 31
$ws.Range('E31').Value = '  -4.50%  '

# Row 32
$ws.Range('E32').Value = '  -1.63%  '

# Row 33
$ws.Range('E33').Value = '  -2.86%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.386.06'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.24%  '

# Row 35
$ws.Range('E35').Value = '  +0.81%  '

# Row 36
$ws.Range('E36').Value = '  -3.44%  '

# Row 37
$ws.Range('E37').Value = '  -1.35%  '

# Row 38
$ws.Range('E38').Value = '  -1.14%  '

# Row 39
$ws.Range('E39').Value = '  -3.18%  '

# Row 40
$ws.Range('E40').Value = '  +0.55%  '

# Row 41
$ws.Range('E41').Value = '  -3.33%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.999'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.12%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.778'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.38%  '

# Row 44
$ws.Range('E44').Value = '  -0.66%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '5.39'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.41%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '61.93'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.26%  '

# Row 47
$ws.Range('B47').Value = 'WEMIXToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.907'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -6.19%  '

# Row 48
$ws.Range('B48').Value = 'RocketPoolETH'
$ws.Range('C48').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.687.44'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.41%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '85.68'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.24%  '

# Row 50
$ws.Range('B50').Value = 'BitcoinSV'
$ws.Range('C50').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '42.12'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +5.24%  '

# Row 51
$ws.Range('B51').Value = 'BabyDogeCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0₆0103'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.64%  '
